# Generate Report for Handoff
# - Mark the newly-handed-off files (rows 7,9,10,12,13,14) with the "ht"
#   (handoff type) priority flag on both locale sheets.
# - Refresh the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#   timestamps for those same rows to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 12, 13, 14)

# zh-cn sheet: Priority (E) -> "ht", Latest Handoff Datetime (H) -> new time
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-03 12:23:28"
}

# de-de sheet: Priority (E) -> "ht", Latest Handoff Datetime (H) -> new time
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-03 12:23:33"
}

# Overview sheet: Latest HO Xliff Generate Date (G) -> new time (matches de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-03 12:23:33"
}
